$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SRLV_curated.csv")

# Rows 676-698: full records for new Ovis/USA genome sequences (N. America region)
$accessionsFull = @(
    "MT993918",
    "MT993917",
    "MT993916",
    "MT993915",
    "MT993914",
    "MT993913",
    "MT993912",
    "MT993911",
    "MT993910",
    "MT993909",
    "MT993908",
    "MT993907",
    "MT993906",
    "MT993905",
    "MT993904",
    "MT993903",
    "MT993902",
    "MT993901",
    "MT993900",
    "MT993899",
    "MT993898",
    "MT993897",
    "MT993896"
)

$startRowFull = 676
for ($i = 0; $i -lt $accessionsFull.Length; $i++) {
    $r = $startRowFull + $i
    $ws.Cells.Item($r, 2).Value = $accessionsFull[$i]   # B: Accession
    $ws.Cells.Item($r, 3).Value = "Ovis"                # C: Species
    $ws.Cells.Item($r, 4).Value = "USA"                 # D: Country
    $ws.Cells.Item($r, 7).Value = "A"                   # G: Genotype
    $ws.Cells.Item($r, 11).Value = "Genome"             # K: Genomic region
    $ws.Cells.Item($r, 12).Value = "N. America"         # L: Region
}

# Rows 699-723: partial records (species/genomic region only) for newly incorporated sequences
$accessionsPartial = @(
    "MG554414",
    "MG554413",
    "MG554412",
    "MG554411",
    "MG554410",
    "MG554409",
    "MG554408",
    "MG554407",
    "MG554406",
    "MG554405",
    "MG554404",
    "MG554403",
    "MG554402",
    "MH374291",
    "MH374290",
    "MH374289",
    "MH374288",
    "MH374287",
    "MH374286",
    "MH374285",
    "MH374284",
    "MH374283",
    "MH936675",
    "MH936674",
    "MG996440"
)

$startRowPartial = 699
for ($i = 0; $i -lt $accessionsPartial.Length; $i++) {
    $r = $startRowPartial + $i
    $ws.Cells.Item($r, 2).Value = $accessionsPartial[$i]  # B: Accession
    $ws.Cells.Item($r, 3).Value = "NK"                    # C: Species
    $ws.Cells.Item($r, 11).Value = "Genome"               # K: Genomic region
}

# Column B (Accession) cells were pasted without the column's default left-aligned
# style, unlike the other populated columns in these new rows - match that.
$ws.Range("B676:B723").Style = "Normal"

# Reposition the view as left by the author after the edit
$win = $excel.ActiveWindow
$win.ScrollRow = 671
$win.ScrollColumn = 1
[void]$ws.Range("I697").Select()
